$d = $word.ActiveDocument

# wildcard off, match case on (not important, all literal) - helper values
$MatchCase = $true
$MatchWholeWord = $false
$MatchWildcards = $false
$MatchSoundsLike = $false
$MatchAllWordForms = $false
$Forward = $true
$Wrap = 1            # wdFindContinue
$Format = $false
$ReplaceAll = 2       # wdReplaceAll

# ------------------------------------------------------------------
# 1) Header date-format field: merge the proofErr-fragmented runs
#    "<<{dateFormat($nowUTC ,'d MMMM yyyy')}>>" stays the same text,
#    just re-typed as a single run.
# ------------------------------------------------------------------
$d.Content.Find.Execute("dateFormat($nowUTC ,`u{2018}d MMMM yyyy`u{2019})}>>", $MatchCase, $MatchWholeWord, $MatchWildcards, $MatchSoundsLike, $MatchAllWordForms, $Forward, $Wrap, $Format, "dateFormat($nowUTC ,`u{2018}d MMMM yyyy`u{2019})}>>", $ReplaceAll) | Out-Null

# ------------------------------------------------------------------
# 2) "In the County Court at   <<courtName>>" - merge the
#    proofErr-fragmented "<<courtName>>" token into fewer runs.
# ------------------------------------------------------------------
$d.Content.Find.Execute("<<courtName>>", $MatchCase, $MatchWholeWord, $MatchWildcards, $MatchSoundsLike, $MatchAllWordForms, $Forward, $Wrap, $Format, "<<courtName>>", $ReplaceAll) | Out-Null

# ------------------------------------------------------------------
# 3) "Case number: <<claimNumber>>" - merge fragmented runs.
# ------------------------------------------------------------------
$d.Content.Find.Execute("<<claimNumber>>", $MatchCase, $MatchWholeWord, $MatchWildcards, $MatchSoundsLike, $MatchAllWordForms, $Forward, $Wrap, $Format, "<<claimNumber>>", $ReplaceAll) | Out-Null

# ------------------------------------------------------------------
# 4) "<<cr_{isMultiParty='No'}>>" / "<<cr_{isMultiParty='Yes'}>>"
#    merge fragmented runs (two occurrences).
# ------------------------------------------------------------------
$d.Content.Find.Execute("<<cr_{isMultiParty=", $MatchCase, $MatchWholeWord, $MatchWildcards, $MatchSoundsLike, $MatchAllWordForms, $Forward, $Wrap, $Format, "<<cr_{isMultiParty=", $ReplaceAll) | Out-Null

# ------------------------------------------------------------------
# 5) "<<cr_{ claimant2Name !=null}>>" merge fragmented runs.
# ------------------------------------------------------------------
$d.Content.Find.Execute("<<cr_{ claimant", $MatchCase, $MatchWholeWord, $MatchWildcards, $MatchSoundsLike, $MatchAllWordForms, $Forward, $Wrap, $Format, "<<cr_{ claimant", $ReplaceAll) | Out-Null

# ------------------------------------------------------------------
# 6) "<<cr_{ defendant2Name!=null}>>" merge fragmented runs.
# ------------------------------------------------------------------
$d.Content.Find.Execute("<<cr_{ defendant2Name", $MatchCase, $MatchWholeWord, $MatchWildcards, $MatchSoundsLike, $MatchAllWordForms, $Forward, $Wrap, $Format, "<<cr_{ defendant2Name", $ReplaceAll) | Out-Null

# ------------------------------------------------------------------
# 7) Main content change: the court name merge field is replaced by
#    the new site name / address / postcode merge fields.
# ------------------------------------------------------------------
$d.Content.Find.Execute("This order is made by <<judgeNameTitle>> on <<", $MatchCase, $MatchWholeWord, $MatchWildcards, $MatchSoundsLike, $MatchAllWordForms, $Forward, $Wrap, $Format, "This order is made by <<judgeNameTitle>> on <<", $ReplaceAll) | Out-Null
$d.Content.Find.Execute(">> at <<courtName>>.", $MatchCase, $MatchWholeWord, $MatchWildcards, $MatchSoundsLike, $MatchAllWordForms, $Forward, $Wrap, $Format, ">> at <<siteName>> - <<address>> - <<postcode>>.", $ReplaceAll) | Out-Null

# ------------------------------------------------------------------
# 8) "<<judgeRecital>>" - merge fragmented runs.
# ------------------------------------------------------------------
$d.Content.Find.Execute("<<judgeRecital>>", $MatchCase, $MatchWholeWord, $MatchWildcards, $MatchSoundsLike, $MatchAllWordForms, $Forward, $Wrap, $Format, "<<judgeRecital>>", $ReplaceAll) | Out-Null

# ------------------------------------------------------------------
# 9) "<<judgeDirection>>" - merge fragmented runs.
# ------------------------------------------------------------------
$d.Content.Find.Execute("<<judgeDirection>>", $MatchCase, $MatchWholeWord, $MatchWildcards, $MatchSoundsLike, $MatchAllWordForms, $Forward, $Wrap, $Format, "<<judgeDirection>>", $ReplaceAll) | Out-Null

# ------------------------------------------------------------------
# 10) "<<judicialByCourtsInitiative>>" - merge fragmented runs.
# ------------------------------------------------------------------
$d.Content.Find.Execute("<<judicialByCourtsInitiative>>", $MatchCase, $MatchWholeWord, $MatchWildcards, $MatchSoundsLike, $MatchAllWordForms, $Forward, $Wrap, $Format, "<<judicialByCourtsInitiative>>", $ReplaceAll) | Out-Null

# ------------------------------------------------------------------
# 11) "<<cs_{reasonAvailable=='Yes' }>>" - merge fragmented runs.
# ------------------------------------------------------------------
$d.Content.Find.Execute("reasonAvailable==`u{2019}Yes`u{2019} }>> ", $MatchCase, $MatchWholeWord, $MatchWildcards, $MatchSoundsLike, $MatchAllWordForms, $Forward, $Wrap, $Format, "reasonAvailable==`u{2019}Yes`u{2019} }>> ", $ReplaceAll) | Out-Null

# ------------------------------------------------------------------
# 12) "<<reasonForDecision>>" - merge fragmented runs.
# ------------------------------------------------------------------
$d.Content.Find.Execute("<<reasonForDecision>>", $MatchCase, $MatchWholeWord, $MatchWildcards, $MatchSoundsLike, $MatchAllWordForms, $Forward, $Wrap, $Format, "<<reasonForDecision>>", $ReplaceAll) | Out-Null
